$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.925.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.65%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.351.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.99%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.88%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.671"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.21%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "72.62"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.95%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.592"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.50%  "
# Row 10
$ws.Range("E10").Value = "  -3.78%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "58.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.71%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.47%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.55%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.702.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.53%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.62%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.906"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.94%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.353.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.65%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.833.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.94%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.72%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.32%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "78.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.63%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.37%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.28%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.73"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.05%  "
# Row 27
$ws.Range("E27").Value = "  -2.59%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.97%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.04%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.85%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.42"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.43%  "
# Row 32
$ws.Range("E32").Value = "  -3.36%  "
# Row 33
$ws.Range("E33").Value = "  -0.87%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0743"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.15%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.12%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.80%  "
# Row 37
$ws.Range("E37").Value = "  -4.05%  "
# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.31%  "
# Row 39
$ws.Range("B39").Value = "THORChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.41"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.48%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0274"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.60%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "65.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.08%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.09%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.32%  "
# Row 44
$ws.Range("E44").Value = "  +5.64%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.199"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.33%  "
# Row 47
$ws.Range("E47").Value = "  -0.14%  "
# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.43%  "
# Row 49
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.12%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.31%  "
# Row 51
$ws.Range("E51").Value = "  -6.80%  "
